{"js": "// Adds the \"Texto do arquivo doc estruturado\" content: a short note\n// paragraph followed by three sections (Introdu\u00e7\u00e3o, Material e M\u00e9todos,\n// Dados do Solo) with bookmarked Heading1 titles, mirroring the OOXML\n// produced by the pandoc/bookdown docx writer for this commit.\n\nconst body = context.document.body;\n\n// Small helper: repeat a sentence N times joined by single spaces, just\n// like the source paragraphs (\"Texto. Texto. Texto. Texto. Texto. Texto.\"\n// repeated and separated by plain-space runs).\nfunction repeatSentence(times) {\n  const sentence = \"Texto. Texto. Texto. Texto. Texto. Texto.\";\n  return new Array(times).fill(sentence).join(\" \");\n}\n\n// 1) Plain note paragraph right after the date line.\nconst pNote = body.insertParagraph(\n  \"Mover o arquivo word dentro da pasta docs\",\n  Word.InsertLocation.end\n);\npNote.style = \"FirstParagraph\";\nawait context.sync();\n\n// Helper to append a Heading1 paragraph with a zero-length bookmark\n// (bookmarkStart/bookmarkEnd pair) anchored at its very start, then the\n// heading text \u2014 this matches <w:bookmarkStart/><w:bookmarkEnd/><w:r>\u2026\nasync function insertBookmarkedHeading(text, bookmarkName) {\n  body.insertParagraph(text, Word.InsertLocation.end).style = \"Heading1\";\n  await context.sync();\n\n  // Re-fetch the paragraph from the collection so the range we collapse\n  // reflects the committed content (a freshly-returned insertParagraph\n  // proxy can still carry a \"whole new paragraph\" span).\n  const paras = body.paragraphs;\n  paras.load(\"items\");\n  await context.sync();\n  const headingPara = paras.items[paras.items.length - 1];\n\n  const startRange = headingPara.getRange(Word.RangeLocation.start);\n  startRange.insertBookmark(bookmarkName);\n  await context.sync();\n}\n\n// 2) \"Introdu\u00e7\u00e3o\" heading + three body paragraphs.\nawait insertBookmarkedHeading(\"Introdu\u00e7\u00e3o\", \"introducao\");\n\nbody.insertParagraph(\"Texto de exemplo\", Word.InsertLocation.end).style = \"FirstParagraph\";\nbody.insertParagraph(\"Texto de exemplo\", Word.InsertLocation.end).style = \"BodyText\";\nbody.insertParagraph(\"Texto de exemplo\", Word.InsertLocation.end).style = \"BodyText\";\nbody.insertParagraph(\"Texto de exemplo\", Word.InsertLocation.end).style = \"BodyText\";\nawait context.sync();\n\n// 3) \"Material e M\u00e9todos\" heading + one long body paragraph (9 repeats).\nawait insertBookmarkedHeading(\"Material e M\u00e9todos\", \"material-e-metodos\");\n\nbody.insertParagraph(repeatSentence(9), Word.InsertLocation.end).style = \"FirstParagraph\";\nawait context.sync();\n\n// 4) \"Dados do Solo\" heading + one long body paragraph (4 repeats).\nawait insertBookmarkedHeading(\"Dados do Solo\", \"dados-do-solo\");\n\nbody.insertParagraph(repeatSentence(4), Word.InsertLocation.end).style = \"FirstParagraph\";\nawait context.sync();\n", "ps1": "# Adds the \"Texto do arquivo doc estruturado\" content: a short note\n# paragraph followed by three sections (Introducao, Material e Metodos,\n# Dados do Solo) with bookmarked Heading1 titles, mirroring the OOXML\n# produced by the pandoc/bookdown docx writer for this commit.\n\n$d = $word.ActiveDocument\n$sel = $word.Selection\n\nfunction Add-Paragraph($text, $style) {\n    $sel.EndKey(6)\n    $sel.TypeParagraph()\n    $sel.Style = $style\n    $sel.TypeText($text)\n}\n\nfunction Add-BookmarkedHeading($text, $bookmarkName) {\n    Add-Paragraph $text \"Heading1\"\n    $p = $d.Paragraphs.Last\n    $r = $p.Range.Duplicate\n    $r.Collapse(1)\n    $d.Bookmarks.Add($bookmarkName, $r)\n}\n\nfunction Repeat-Sentence($times) {\n    $sentence = \"Texto. Texto. Texto. Texto. Texto. Texto.\"\n    $parts = @()\n    for ($i = 0; $i -lt $times; $i++) {\n        $parts += $sentence\n    }\n    return [string]::Join(\" \", $parts)\n}\n\n# 1) Plain note paragraph right after the date line.\nAdd-Paragraph \"Mover o arquivo word dentro da pasta docs\" \"FirstParagraph\"\n\n# 2) \"Introducao\" heading + three body paragraphs.\nAdd-BookmarkedHeading \"Introdu\u00e7\u00e3o\" \"introducao\"\nAdd-Paragraph \"Texto de exemplo\" \"FirstParagraph\"\nAdd-Paragraph \"Texto de exemplo\" \"BodyText\"\nAdd-Paragraph \"Texto de exemplo\" \"BodyText\"\nAdd-Paragraph \"Texto de exemplo\" \"BodyText\"\n\n# 3) \"Material e Metodos\" heading + one long body paragraph (9 repeats).\nAdd-BookmarkedHeading \"Material e M\u00e9todos\" \"material-e-metodos\"\nAdd-Paragraph (Repeat-Sentence 9) \"FirstParagraph\"\n\n# 4) \"Dados do Solo\" heading + one long body paragraph (4 repeats).\nAdd-BookmarkedHeading \"Dados do Solo\" \"dados-do-solo\"\nAdd-Paragraph (Repeat-Sentence 4) \"FirstParagraph\"\n"}
